# Applies the Unicorn_Profits.xlsx scheduled-runner value refresh
# (commit: "chore: update Sheets via scheduled runner") across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
#
# Every touched cell is a plain numeric literal (no formulas in this
# workbook), so each target is written with a direct Range.Value
# assignment; the single cell that the refresh actually blanks out
# (WVR!N56) is cleared with ClearContents so the <c> element drops
# out of the sheet XML entirely, matching the source diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 166667200
$ws.Range("I8").Value = 200000030
$ws.Range("K8").Value = 600000090
$ws.Range("M8").Value = -599999951
$ws.Range("H17").Value = 799.0411
$ws.Range("J17").Value = 785.4219000000001
$ws.Range("L17").Value = 2356.2657
$ws.Range("N17").Value = -2692.2657
$ws.Range("H112").Value = 1717.4615
$ws.Range("J112").Value = 1754.16
$ws.Range("L112").Value = 5262.48
$ws.Range("N112").Value = -7478.48
$ws.Range("H127").Value = 645.4
$ws.Range("I127").Value = 265
$ws.Range("J127").Value = 1216
$ws.Range("K127").Value = 795
$ws.Range("L127").Value = 3648
$ws.Range("M127").Value = 4165
$ws.Range("N127").Value = -13568

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 83.333336
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -274
$ws.Range("H11").Value = 16500
$ws.Range("I11").Value = 8000
$ws.Range("J11").Value = 25000
$ws.Range("K11").Value = 8000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = -7856
$ws.Range("N11").Value = -25288
$ws.Range("H61").Value = 297683.4
$ws.Range("I61").Value = 251327.1
$ws.Range("J61").Value = 373227.03
$ws.Range("K61").Value = 251327.1
$ws.Range("L61").Value = 373227.03
$ws.Range("M61").Value = -251115.1
$ws.Range("N61").Value = -373651.03
$ws.Range("H74").Value = 144272.92
$ws.Range("I74").Value = 150109.31
$ws.Range("J74").Value = 105169.1
$ws.Range("K74").Value = 150109.31
$ws.Range("L74").Value = 105169.1
$ws.Range("M74").Value = -149235.31
$ws.Range("N74").Value = -106917.1
$ws.Range("H77").Value = 144272.92
$ws.Range("I77").Value = 150109.31
$ws.Range("J77").Value = 105169.1
$ws.Range("K77").Value = 750546.55
$ws.Range("L77").Value = 525845.5
$ws.Range("M77").Value = -746178.55
$ws.Range("N77").Value = -534581.5
$ws.Range("H132").Value = 15301.507
$ws.Range("I132").Value = 19091.965
$ws.Range("K132").Value = 57275.895
$ws.Range("M132").Value = -54745.895
$ws.Range("H136").Value = 297683.4
$ws.Range("I136").Value = 251327.1
$ws.Range("J136").Value = 373227.03
$ws.Range("K136").Value = 753981.3
$ws.Range("L136").Value = 1119681.09
$ws.Range("M136").Value = -751431.3
$ws.Range("N136").Value = -1124781.09

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 83.333336
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -280

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 800
$ws.Range("I12").Value = 800
$ws.Range("K12").Value = 800
$ws.Range("M12").Value = -630
$ws.Range("H31").Value = 2091.423
$ws.Range("I31").Value = 1489.8813
$ws.Range("K31").Value = 1489.8813
$ws.Range("M31").Value = -1194.8813
$ws.Range("H34").Value = 2091.423
$ws.Range("I34").Value = 1489.8813
$ws.Range("K34").Value = 1489.8813
$ws.Range("M34").Value = -1287.8813
$ws.Range("H58").Value = 3495.8333
$ws.Range("I58").Value = 3715.4849
$ws.Range("K58").Value = 3715.4849
$ws.Range("M58").Value = -3512.4849
$ws.Range("H92").Value = 29560
$ws.Range("J92").Value = 29560
$ws.Range("L92").Value = 29560
$ws.Range("N92").Value = -34552
$ws.Range("H99").Value = 127297.25
$ws.Range("I99").Value = 167918.67
$ws.Range("J99").Value = 5433
$ws.Range("K99").Value = 167918.67
$ws.Range("L99").Value = 5433
$ws.Range("M99").Value = -166420.67
$ws.Range("N99").Value = -8429
$ws.Range("H126").Value = 127297.25
$ws.Range("I126").Value = 167918.67
$ws.Range("J126").Value = 5433
$ws.Range("K126").Value = 503756.01
$ws.Range("L126").Value = 16299
$ws.Range("M126").Value = -501286.01
$ws.Range("N126").Value = -21239
$ws.Range("H134").Value = 1166.197
$ws.Range("I134").Value = 750.4902
$ws.Range("J134").Value = 2579.6
$ws.Range("K134").Value = 2251.4706
$ws.Range("L134").Value = 7738.799999999999
$ws.Range("M134").Value = 283.5294000000004
$ws.Range("N134").Value = -12808.8
$ws.Range("H136").Value = 3495.8333
$ws.Range("I136").Value = 3715.4849
$ws.Range("K136").Value = 11146.4547
$ws.Range("M136").Value = -8596.4547

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 498.42856
$ws.Range("I5").Value = 468.35
$ws.Range("K5").Value = 1405.05
$ws.Range("M5").Value = -1293.05
$ws.Range("H135").Value = 498.42856
$ws.Range("I135").Value = 468.35
$ws.Range("K135").Value = 4215.150000000001
$ws.Range("M135").Value = -1680.150000000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4565.552
$ws.Range("I113").Value = 5805.6
$ws.Range("J113").Value = 1809.8889
$ws.Range("K113").Value = 5805.6
$ws.Range("L113").Value = 1809.8889
$ws.Range("M113").Value = -3635.6
$ws.Range("N113").Value = -6149.8889
$ws.Range("H132").Value = 3501.8408
$ws.Range("I132").Value = 3835.1304
$ws.Range("J132").Value = 3136.8096
$ws.Range("K132").Value = 11505.3912
$ws.Range("L132").Value = 9410.4288
$ws.Range("M132").Value = -8975.3912
$ws.Range("N132").Value = -14470.4288

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7253.372
$ws.Range("I132").Value = 2418.5217
$ws.Range("J132").Value = 12813.45
$ws.Range("K132").Value = 7255.5651
$ws.Range("L132").Value = 38440.35000000001
$ws.Range("M132").Value = -4725.5651
$ws.Range("N132").Value = -43500.35000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 4950
$ws.Range("J18").Value = 3666.6667
$ws.Range("L18").Value = 3666.6667
$ws.Range("N18").Value = -4012.6667
$ws.Range("H56").Value = 3900
$ws.Range("I56").Value = 3900
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 3900
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -3186
$ws.Range("N56").ClearContents()  # was -15053, now blank
$ws.Range("H64").Value = 9900
$ws.Range("J64").Value = 9900
$ws.Range("L64").Value = 9900
$ws.Range("N64").Value = -10396
$ws.Range("H67").Value = 9900
$ws.Range("J67").Value = 9900
$ws.Range("L67").Value = 9900
$ws.Range("N67").Value = -11616
$ws.Range("H126").Value = 1692.3684
$ws.Range("I126").Value = 763.5454999999999
$ws.Range("J126").Value = 2969.5
$ws.Range("K126").Value = 2290.6365
$ws.Range("L126").Value = 8908.5
$ws.Range("M126").Value = 179.3635000000004
$ws.Range("N126").Value = -13848.5
$ws.Range("H132").Value = 1456.5692
$ws.Range("I132").Value = 788.6539
$ws.Range("J132").Value = 4128.231
$ws.Range("K132").Value = 2365.9617
$ws.Range("L132").Value = 12384.693
$ws.Range("M132").Value = 164.0383000000002
$ws.Range("N132").Value = -17444.693

Write-Output "Updated 178 cells and cleared 1 cell across 8 sheets."
